$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.728.85"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "1.596.65"
$ws.Range("E3").Value = "  -1.63%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.60"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("E6").Value = "  -2.41%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -2.56%  "
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("E10").Value = "  -1.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0869"
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("D12").Value = "1.822.85"
$ws.Range("E12").Value = "  -1.71%  "
$ws.Range("D13").Value = "1.591.67"
$ws.Range("E13").Value = "  -2.25%  "
$ws.Range("E14").Value = "  -3.28%  "
$ws.Range("E15").Value = "  -3.03%  "
$ws.Range("D16").Value = "27.729.78"
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.59"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "219.81"
$ws.Range("E18").Value = "  -3.38%  "
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.38"
$ws.Range("E20").Value = "  -2.85%  "
$ws.Range("E22").Value = "  -3.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.74"
$ws.Range("E23").Value = "  -1.96%  "
$ws.Range("E24").Value = "  -3.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.22"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.78"
$ws.Range("E26").Value = "  -1.71%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  -1.13%  "
$ws.Range("E29").Value = "  -4.43%  "
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("E31").Value = "  -1.47%  "
$ws.Range("E32").Value = "  -4.53%  "
$ws.Range("D33").Value = "1.379.05"
$ws.Range("E33").Value = "  -2.38%  "
$ws.Range("E34").Value = "  -3.58%  "
$ws.Range("E35").Value = "  -3.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.972"
$ws.Range("E36").Value = "  -2.08%  "
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.537"
$ws.Range("E39").Value = "  -2.77%  "
$ws.Range("E40").Value = "  -2.06%  "
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("E42").Value = "  -2.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.68"
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("E44").Value = "  +2.74%  "
$ws.Range("E45").Value = "  -2.95%  "
$ws.Range("E46").Value = "  -4.11%  "
$ws.Range("D47").Value = "1.733.32"
$ws.Range("E47").Value = "  -1.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.06"
$ws.Range("E48").Value = "  -2.24%  "
$ws.Range("D49").Value = "0.0₆0102"
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("E50").Value = "  -3.48%  "
$ws.Range("E51").Value = "  -1.22%  "
